$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 655.375
$ws.Range("I28").Value = 605.2857
$ws.Range("K28").Value = 605.2857
$ws.Range("M28").Value = -120.2857
$ws.Range("H33").Value = 224.4
$ws.Range("I33").Value = 224.4
$ws.Range("K33").Value = 224.4
$ws.Range("M33").Value = 4.599999999999994
$ws.Range("H93").Value = 25000
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -29992
$ws.Range("H96").Value = 2094
$ws.Range("I96").Value = 3090
$ws.Range("J96").Value = 600
$ws.Range("K96").Value = 9270
$ws.Range("L96").Value = 1800
$ws.Range("M96").Value = -7897
$ws.Range("N96").Value = -4546
$ws.Range("H98").Value = 1148.1818
$ws.Range("I98").Value = 1148.1818
$ws.Range("K98").Value = 1148.1818
$ws.Range("M98").Value = 349.8181999999999
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 3800
$ws.Range("J113").Value = 3800
$ws.Range("L113").Value = 3800
$ws.Range("N113").Value = -10308
$ws.Range("H122").Value = 1148.1818
$ws.Range("I122").Value = 1148.1818
$ws.Range("K122").Value = 3444.5454
$ws.Range("M122").Value = -994.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 1104.25
$ws.Range("J17").Value = 2008.5
$ws.Range("L17").Value = 2008.5
$ws.Range("N17").Value = -2354.5
$ws.Range("H61").Value = 3071.1428
$ws.Range("I61").Value = 3071.1428
$ws.Range("K61").Value = 3071.1428
$ws.Range("M61").Value = -2859.1428
$ws.Range("H74").Value = 3885.4285
$ws.Range("I74").Value = 2459.6
$ws.Range("K74").Value = 2459.6
$ws.Range("M74").Value = -1585.6
$ws.Range("H77").Value = 3885.4285
$ws.Range("I77").Value = 2459.6
$ws.Range("K77").Value = 12298
$ws.Range("M77").Value = -7930
$ws.Range("H97").Value = 3373.75
$ws.Range("I97").Value = 1898.3334
$ws.Range("K97").Value = 1898.3334
$ws.Range("M97").Value = -1402.3334
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H110").Value = 996
$ws.Range("I110").Value = 994
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 994
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1051
$ws.Range("N110").Value = -5090
$ws.Range("H122").Value = 2666
$ws.Range("J122").Value = 2499
$ws.Range("L122").Value = 7497
$ws.Range("N122").Value = -12397
$ws.Range("H136").Value = 3071.1428
$ws.Range("I136").Value = 3071.1428
$ws.Range("K136").Value = 9213.428400000001
$ws.Range("M136").Value = -6663.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 350
$ws.Range("I10").Value = 350
$ws.Range("K10").Value = 350
$ws.Range("M10").Value = -210
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 5000
$ws.Range("K20").Value = 5000
$ws.Range("M20").Value = -4753

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 751.8889
$ws.Range("I16").Value = 720.875
$ws.Range("K16").Value = 720.875
$ws.Range("M16").Value = -433.875
$ws.Range("H22").Value = 280
$ws.Range("H31").Value = 6410.0713
$ws.Range("I31").Value = 3896.6
$ws.Range("J31").Value = 7806.4443
$ws.Range("K31").Value = 3896.6
$ws.Range("L31").Value = 7806.4443
$ws.Range("M31").Value = -3601.6
$ws.Range("N31").Value = -8396.444299999999
$ws.Range("H34").Value = 6410.0713
$ws.Range("I34").Value = 3896.6
$ws.Range("J34").Value = 7806.4443
$ws.Range("K34").Value = 3896.6
$ws.Range("L34").Value = 7806.4443
$ws.Range("M34").Value = -3694.6
$ws.Range("N34").Value = -8210.444299999999
$ws.Range("H41").Value = 15082.5
$ws.Range("J41").Value = 20165
$ws.Range("L41").Value = 20165
$ws.Range("N41").Value = -21021
$ws.Range("H51").Value = 24333.334
$ws.Range("J51").Value = 24333.334
$ws.Range("L51").Value = 24333.334
$ws.Range("N51").Value = -25805.334
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 24333.334
$ws.Range("J61").Value = 24333.334
$ws.Range("L61").Value = 24333.334
$ws.Range("N61").Value = -25029.334
$ws.Range("H107").Value = 625
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 751.8889
$ws.Range("I113").Value = 720.875
$ws.Range("K113").Value = 720.875
$ws.Range("M113").Value = 1449.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 792.6667
$ws.Range("J113").Value = 600
$ws.Range("L113").Value = 1800
$ws.Range("N113").Value = -6140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4999
$ws.Range("J80").Value = 4999
$ws.Range("L80").Value = 4999
$ws.Range("N80").Value = -6995
$ws.Range("H83").Value = 4999
$ws.Range("J83").Value = 4999
$ws.Range("L83").Value = 24995
$ws.Range("N83").Value = -34979
$ws.Range("H97").Value = 1539
$ws.Range("J97").Value = 691.6667
$ws.Range("L97").Value = 691.6667
$ws.Range("N97").Value = -1683.6667
$ws.Range("H122").Value = 2899.8572
$ws.Range("I122").Value = 1299.75
$ws.Range("K122").Value = 3899.25
$ws.Range("M122").Value = -1449.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 14450
$ws.Range("I33").Value = 100
$ws.Range("K33").Value = 100
$ws.Range("M33").Value = 190
$ws.Range("H40").Value = 8465
$ws.Range("I40").Value = 8472.5
$ws.Range("J40").Value = 8450
$ws.Range("K40").Value = 8472.5
$ws.Range("L40").Value = 8450
$ws.Range("M40").Value = -8336.5
$ws.Range("N40").Value = -8722
$ws.Range("H55").Value = 1498.7
$ws.Range("J55").Value = 2532.3333
$ws.Range("L55").Value = 2532.3333
$ws.Range("N55").Value = -2878.3333
$ws.Range("H97").Value = 8844
$ws.Range("J97").Value = 8844
$ws.Range("L97").Value = 8844
$ws.Range("N97").Value = -10826
$ws.Range("H100").Value = 16250
$ws.Range("I100").Value = 20000
$ws.Range("K100").Value = 20000
$ws.Range("M100").Value = -19459
$ws.Range("H122").Value = 5064.8335
$ws.Range("J122").Value = 3996.3333
$ws.Range("L122").Value = 11988.9999
$ws.Range("N122").Value = -16888.9999
$ws.Range("H125").Value = 99995
$ws.Range("J125").Value = 99995
$ws.Range("L125").Value = 99995
$ws.Range("N125").Value = -109835
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23333.334
$ws.Range("I54").Value = 10000
$ws.Range("K54").Value = 10000
$ws.Range("M54").Value = -9480
$ws.Range("H107").Value = 2217.4443
$ws.Range("J107").Value = 4601.5
$ws.Range("L107").Value = 13804.5
$ws.Range("N107").Value = -17644.5
$ws.Range("H126").Value = 500
$ws.Range("I126").Value = 500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 970
$ws.Range("N126").ClearContents()
